$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743637847690/PaZ2GKc6aN-camera-video-segment-driver-1743637850190.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE43XEL5%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T200045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAvE7mpeXw2hL46foMlQaVZu%2FrDL%2FwIAdifU%2BmwyYnXGAiEAsvrsJpw35QlrxguEQsYycT7lUn3DBL183LmXRKt7c60q5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLrn%2B8y9b%2FjaDK41eyq6A4MuNEZbmp4DVKmjBfgtyXdzGQ%2FjdZ6%2Fg6eDLy9QktJ4TC8msZw8RPftWfZsRsk9tibfy303ciELIipxbCpY1w86Trccc15%2BUuISfPSzFWb8PHGKj2I8x9c7VqCeBccYEo5ItCqoQd7MgYsG5cVH16JqY9n5PWwUCUvmn6iCRke2M1oDJQDThwiruZ7FFB7Y7kqusTLxFjgMI3lqODycybRZIehp%2F%2FAkTu%2Be7FAQb130i5aSHdsBGp0Eijh6fgA9I6cpx%2BSOcFnuy3UiZC38Yvcnn468zL2493K3DjSfnlYKwr6ufZP48pgA8WgH6ABmrC0MmEv5o%2FgEL7RGF%2Ff%2B6%2BKsovaAeLgamULS3ixAKil5UrOt2iMSv9YxuU80udkXlTHknzdzjPTA9J2I25HFMQJcnrzFOkMjGGstSEZLxVqJ5HgT5CMP23N8gM7BbpGpj%2BIxtCBHXxN7pbdhYd3JkzdVsJ4UbKBiQj%2F7zG3cfokxuQR8DwLBp4Fxu5jM%2BmK5ORlebTmuGxb9BsXwQDfUHMuNdH4CWCmGL0KVOeGHzdMdDT4OWxoLNkdXPVgAsev9GCM%2BlP%2BLUSudhYgwy8O7vwY6pQGtk40YKaz4nFRgL4UDORQGxF53O%2Fs9bHVl%2FNsnNTLCX5ySvJlrT67jHLUW6jR2o7DBf8coqoxmdBbXuSReCAnAht7H8Bh%2FsW%2BeBsvNoLipmBnn8BzO3ty4ctWHHHsvrQD%2BIU%2FWpwoNeUoYKPEA3MeDFu06U58heS3XRO9lDWqKjJR9LxEgE5isnCQiQtQNnw%2Bj1AemosTzMXCV7Vj%2BDFOCnH1iREo%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2004%3A00%3A45%20GMT&X-Amz-Signature=c2c298cdb7f41fc95ea371f9d91ae50ec0cd85ebd47ecdfc455c3bc028bd7d98'

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743634721601/5dyCDJLTC3-camera-video-segment-driver-1743634724101.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE43XEL5%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T200045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAvE7mpeXw2hL46foMlQaVZu%2FrDL%2FwIAdifU%2BmwyYnXGAiEAsvrsJpw35QlrxguEQsYycT7lUn3DBL183LmXRKt7c60q5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLrn%2B8y9b%2FjaDK41eyq6A4MuNEZbmp4DVKmjBfgtyXdzGQ%2FjdZ6%2Fg6eDLy9QktJ4TC8msZw8RPftWfZsRsk9tibfy303ciELIipxbCpY1w86Trccc15%2BUuISfPSzFWb8PHGKj2I8x9c7VqCeBccYEo5ItCqoQd7MgYsG5cVH16JqY9n5PWwUCUvmn6iCRke2M1oDJQDThwiruZ7FFB7Y7kqusTLxFjgMI3lqODycybRZIehp%2F%2FAkTu%2Be7FAQb130i5aSHdsBGp0Eijh6fgA9I6cpx%2BSOcFnuy3UiZC38Yvcnn468zL2493K3DjSfnlYKwr6ufZP48pgA8WgH6ABmrC0MmEv5o%2FgEL7RGF%2Ff%2B6%2BKsovaAeLgamULS3ixAKil5UrOt2iMSv9YxuU80udkXlTHknzdzjPTA9J2I25HFMQJcnrzFOkMjGGstSEZLxVqJ5HgT5CMP23N8gM7BbpGpj%2BIxtCBHXxN7pbdhYd3JkzdVsJ4UbKBiQj%2F7zG3cfokxuQR8DwLBp4Fxu5jM%2BmK5ORlebTmuGxb9BsXwQDfUHMuNdH4CWCmGL0KVOeGHzdMdDT4OWxoLNkdXPVgAsev9GCM%2BlP%2BLUSudhYgwy8O7vwY6pQGtk40YKaz4nFRgL4UDORQGxF53O%2Fs9bHVl%2FNsnNTLCX5ySvJlrT67jHLUW6jR2o7DBf8coqoxmdBbXuSReCAnAht7H8Bh%2FsW%2BeBsvNoLipmBnn8BzO3ty4ctWHHHsvrQD%2BIU%2FWpwoNeUoYKPEA3MeDFu06U58heS3XRO9lDWqKjJR9LxEgE5isnCQiQtQNnw%2Bj1AemosTzMXCV7Vj%2BDFOCnH1iREo%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2004%3A00%3A45%20GMT&X-Amz-Signature=255e1daeb47c5edbec7992df8071714086ae309ba5852943d7ac755b590ff7c4'

$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1743618603072/lsfY7Elc8s-camera-video-segment-driver-1743618608072.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE43XEL5%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T200045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAvE7mpeXw2hL46foMlQaVZu%2FrDL%2FwIAdifU%2BmwyYnXGAiEAsvrsJpw35QlrxguEQsYycT7lUn3DBL183LmXRKt7c60q5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLrn%2B8y9b%2FjaDK41eyq6A4MuNEZbmp4DVKmjBfgtyXdzGQ%2FjdZ6%2Fg6eDLy9QktJ4TC8msZw8RPftWfZsRsk9tibfy303ciELIipxbCpY1w86Trccc15%2BUuISfPSzFWb8PHGKj2I8x9c7VqCeBccYEo5ItCqoQd7MgYsG5cVH16JqY9n5PWwUCUvmn6iCRke2M1oDJQDThwiruZ7FFB7Y7kqusTLxFjgMI3lqODycybRZIehp%2F%2FAkTu%2Be7FAQb130i5aSHdsBGp0Eijh6fgA9I6cpx%2BSOcFnuy3UiZC38Yvcnn468zL2493K3DjSfnlYKwr6ufZP48pgA8WgH6ABmrC0MmEv5o%2FgEL7RGF%2Ff%2B6%2BKsovaAeLgamULS3ixAKil5UrOt2iMSv9YxuU80udkXlTHknzdzjPTA9J2I25HFMQJcnrzFOkMjGGstSEZLxVqJ5HgT5CMP23N8gM7BbpGpj%2BIxtCBHXxN7pbdhYd3JkzdVsJ4UbKBiQj%2F7zG3cfokxuQR8DwLBp4Fxu5jM%2BmK5ORlebTmuGxb9BsXwQDfUHMuNdH4CWCmGL0KVOeGHzdMdDT4OWxoLNkdXPVgAsev9GCM%2BlP%2BLUSudhYgwy8O7vwY6pQGtk40YKaz4nFRgL4UDORQGxF53O%2Fs9bHVl%2FNsnNTLCX5ySvJlrT67jHLUW6jR2o7DBf8coqoxmdBbXuSReCAnAht7H8Bh%2FsW%2BeBsvNoLipmBnn8BzO3ty4ctWHHHsvrQD%2BIU%2FWpwoNeUoYKPEA3MeDFu06U58heS3XRO9lDWqKjJR9LxEgE5isnCQiQtQNnw%2Bj1AemosTzMXCV7Vj%2BDFOCnH1iREo%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2004%3A00%3A45%20GMT&X-Amz-Signature=d22cff7206c8f47374c8317ba584a15e91427f4d6739c08079fc91e0c17f314b'

$ws.Range("L5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743618603072/krpw3PxSiC-camera-video-segment-1743618608072.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE43XEL5%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T200045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAvE7mpeXw2hL46foMlQaVZu%2FrDL%2FwIAdifU%2BmwyYnXGAiEAsvrsJpw35QlrxguEQsYycT7lUn3DBL183LmXRKt7c60q5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDLrn%2B8y9b%2FjaDK41eyq6A4MuNEZbmp4DVKmjBfgtyXdzGQ%2FjdZ6%2Fg6eDLy9QktJ4TC8msZw8RPftWfZsRsk9tibfy303ciELIipxbCpY1w86Trccc15%2BUuISfPSzFWb8PHGKj2I8x9c7VqCeBccYEo5ItCqoQd7MgYsG5cVH16JqY9n5PWwUCUvmn6iCRke2M1oDJQDThwiruZ7FFB7Y7kqusTLxFjgMI3lqODycybRZIehp%2F%2FAkTu%2Be7FAQb130i5aSHdsBGp0Eijh6fgA9I6cpx%2BSOcFnuy3UiZC38Yvcnn468zL2493K3DjSfnlYKwr6ufZP48pgA8WgH6ABmrC0MmEv5o%2FgEL7RGF%2Ff%2B6%2BKsovaAeLgamULS3ixAKil5UrOt2iMSv9YxuU80udkXlTHknzdzjPTA9J2I25HFMQJcnrzFOkMjGGstSEZLxVqJ5HgT5CMP23N8gM7BbpGpj%2BIxtCBHXxN7pbdhYd3JkzdVsJ4UbKBiQj%2F7zG3cfokxuQR8DwLBp4Fxu5jM%2BmK5ORlebTmuGxb9BsXwQDfUHMuNdH4CWCmGL0KVOeGHzdMdDT4OWxoLNkdXPVgAsev9GCM%2BlP%2BLUSudhYgwy8O7vwY6pQGtk40YKaz4nFRgL4UDORQGxF53O%2Fs9bHVl%2FNsnNTLCX5ySvJlrT67jHLUW6jR2o7DBf8coqoxmdBbXuSReCAnAht7H8Bh%2FsW%2BeBsvNoLipmBnn8BzO3ty4ctWHHHsvrQD%2BIU%2FWpwoNeUoYKPEA3MeDFu06U58heS3XRO9lDWqKjJR9LxEgE5isnCQiQtQNnw%2Bj1AemosTzMXCV7Vj%2BDFOCnH1iREo%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2004%3A00%3A45%20GMT&X-Amz-Signature=f195e733919b819737af4eebab073e5bf7c5e868d6271702ea32dd4b59b48986'
